$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 2 and 3 swap their Fecha / Volumen / Precio / Unidad / Precio $/Kg / Kg-unidad values
$cols = @("D", "M", "N", "O", "P", "Q", "S", "T")

foreach ($col in $cols) {
    $cell2 = $ws.Range($col + "2")
    $cell3 = $ws.Range($col + "3")
    $tmp = $cell2.Value2
    $cell2.Value = $cell3.Value2
    $cell3.Value = $tmp
}
